$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$orig = $ws.Range("D2").Style
$ws.Range("D2").Value = "'43.242.46"
$ws.Range("D2").Style = $orig
$ws.Range("E2").Value = '  +2.51%  '

$orig = $ws.Range("D3").Style
$ws.Range("D3").Value = "'2.311.22"
$ws.Range("D3").Style = $orig
$ws.Range("E3").Value = '  +1.75%  '

$ws.Range("E4").Value = '  -0.02%  '

$orig = $ws.Range("D5").Style
$ws.Range("D5").Value = "'302.53"
$ws.Range("D5").Style = $orig
$ws.Range("E5").Value = '  +1.35%  '

$orig = $ws.Range("D6").Style
$ws.Range("D6").Value = "'100.04"
$ws.Range("D6").Style = $orig
$ws.Range("E6").Value = '  +5.65%  '

$orig = $ws.Range("D7").Style
$ws.Range("D7").Value = "'0.504"
$ws.Range("D7").Style = $orig
$ws.Range("E7").Value = '  +2.16%  '

$ws.Range("E8").Value = '  -0.06%  '

$orig = $ws.Range("D9").Style
$ws.Range("D9").Value = "'0.509"
$ws.Range("D9").Style = $orig
$ws.Range("E9").Value = '  +3.49%  '

$orig = $ws.Range("D10").Style
$ws.Range("D10").Value = "'34.54"
$ws.Range("D10").Style = $orig
$ws.Range("E10").Value = '  +4.23%  '

$orig = $ws.Range("D11").Style
$ws.Range("D11").Value = "'0.0799"
$ws.Range("D11").Style = $orig
$ws.Range("E11").Value = '  +1.35%  '

$ws.Range("E12").Value = '  +4.69%  '

$orig = $ws.Range("D13").Style
$ws.Range("D13").Value = "'18.08"
$ws.Range("D13").Style = $orig
$ws.Range("E13").Value = '  +15.35%  '

$orig = $ws.Range("D14").Style
$ws.Range("D14").Value = "'6.84"
$ws.Range("D14").Style = $orig
$ws.Range("E14").Value = '  +3.24%  '

$orig = $ws.Range("D15").Style
$ws.Range("D15").Value = "'2.666.92"
$ws.Range("D15").Style = $orig
$ws.Range("E15").Value = '  +1.64%  '

$orig = $ws.Range("D16").Style
$ws.Range("D16").Value = "'2.337.95"
$ws.Range("D16").Style = $orig
$ws.Range("E16").Value = '  +2.17%  '

$orig = $ws.Range("D17").Style
$ws.Range("D17").Value = "'0.819"
$ws.Range("D17").Style = $orig
$ws.Range("E17").Value = '  +5.63%  '

$orig = $ws.Range("D18").Style
$ws.Range("D18").Value = "'43.149.48"
$ws.Range("D18").Style = $orig
$ws.Range("E18").Value = '  +2.32%  '

$orig = $ws.Range("D19").Style
$ws.Range("D19").Value = "'12.69"
$ws.Range("D19").Style = $orig
$ws.Range("E19").Value = '  +11.57%  '

$orig = $ws.Range("D20").Style
$ws.Range("D20").Value = "'0.0₃0908"
$ws.Range("D20").Style = $orig
$ws.Range("E20").Value = '  +2.25%  '

$orig = $ws.Range("D21").Style
$ws.Range("D21").Value = "'6.12"
$ws.Range("D21").Style = $orig
$ws.Range("E21").Value = '  +2.41%  '

$orig = $ws.Range("D22").Style
$ws.Range("D22").Value = "'67.89"
$ws.Range("D22").Style = $orig
$ws.Range("E22").Value = '  +1.92%  '

$orig = $ws.Range("D23").Style
$ws.Range("D23").Value = "'237.77"
$ws.Range("D23").Style = $orig
$ws.Range("E23").Value = '  +2.17%  '

$orig = $ws.Range("D24").Style
$ws.Range("D24").Value = "'2.21"
$ws.Range("D24").Style = $orig
$ws.Range("E24").Value = '  +13.95%  '

$ws.Range("B25").Value = 'PancakeSwap'
$ws.Range("C25").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$orig = $ws.Range("D25").Style
$ws.Range("D25").Value = "'2.47"
$ws.Range("D25").Style = $orig
$ws.Range("E25").Value = '  +0.86%  '

$ws.Range("B26").Value = 'Dai'
$ws.Range("C26").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$orig = $ws.Range("D26").Style
$ws.Range("D26").Value = "'0.999"
$ws.Range("D26").Style = $orig
$ws.Range("E26").Value = '  +0.02%  '

$ws.Range("E27").Value = '  +4.38%  '

$ws.Range("E28").Value = '  -4.90%  '

$orig = $ws.Range("D29").Style
$ws.Range("D29").Value = "'168.64"
$ws.Range("D29").Style = $orig
$ws.Range("E29").Value = '  +1.52%  '

$orig = $ws.Range("D30").Style
$ws.Range("D30").Value = "'34.27"
$ws.Range("D30").Style = $orig
$ws.Range("E30").Value = '  +1.66%  '

$orig = $ws.Range("D31").Style
$ws.Range("D31").Value = "'9.19"
$ws.Range("D31").Style = $orig
$ws.Range("E31").Value = '  +1.72%  '

$ws.Range("E32").Value = '  +0.03%  '

$orig = $ws.Range("D33").Style
$ws.Range("D33").Value = "'5.05"
$ws.Range("D33").Style = $orig
$ws.Range("E33").Value = '  +2.83%  '

$orig = $ws.Range("D34").Style
$ws.Range("D34").Value = "'4.63"
$ws.Range("D34").Style = $orig
$ws.Range("E34").Value = '  +3.76%  '

$orig = $ws.Range("D35").Style
$ws.Range("D35").Value = "'2.44"
$ws.Range("D35").Style = $orig
$ws.Range("E35").Value = '  +4.62%  '

$orig = $ws.Range("D36").Style
$ws.Range("D36").Value = "'17.15"
$ws.Range("D36").Style = $orig
$ws.Range("E36").Value = '  +6.78%  '

$orig = $ws.Range("D37").Style
$ws.Range("D37").Value = "'0.0693"
$ws.Range("D37").Style = $orig
$ws.Range("E37").Value = '  +0.52%  '

$ws.Range("E38").Value = '  +3.88%  '

$orig = $ws.Range("D39").Style
$ws.Range("D39").Value = "'1.81"
$ws.Range("D39").Style = $orig
$ws.Range("E39").Value = '  +5.36%  '

$ws.Range("E40").Value = '  +1.43%  '

$ws.Range("E41").Value = '  +0.77%  '

$orig = $ws.Range("D42").Style
$ws.Range("D42").Value = "'2.37"
$ws.Range("D42").Style = $orig
$ws.Range("E42").Value = '  -2.42%  '

$orig = $ws.Range("D43").Style
$ws.Range("D43").Value = "'2.003.80"
$ws.Range("D43").Style = $orig
$ws.Range("E43").Value = '  +2.47%  '

$orig = $ws.Range("D44").Style
$ws.Range("D44").Value = "'0.0286"
$ws.Range("D44").Style = $orig
$ws.Range("E44").Value = '  +3.42%  '

$orig = $ws.Range("D45").Style
$ws.Range("D45").Value = "'10.15"
$ws.Range("D45").Style = $orig

$orig = $ws.Range("D46").Style
$ws.Range("D46").Value = "'17.75"
$ws.Range("D46").Style = $orig
$ws.Range("E46").Value = '  +1.78%  '

$orig = $ws.Range("D47").Style
$ws.Range("D47").Value = "'2.87"
$ws.Range("D47").Style = $orig
$ws.Range("E47").Value = '  +2.85%  '

$orig = $ws.Range("D48").Style
$ws.Range("D48").Value = "'56.07"
$ws.Range("D48").Style = $orig
$ws.Range("E48").Value = '  +8.01%  '

$orig = $ws.Range("D49").Style
$ws.Range("D49").Value = "'2.535.68"
$ws.Range("D49").Style = $orig
$ws.Range("E49").Value = '  +1.59%  '

$ws.Range("E50").Value = '  +5.29%  '

$orig = $ws.Range("D51").Style
$ws.Range("D51").Value = "'4.58"
$ws.Range("D51").Style = $orig
$ws.Range("E51").Value = '  +1.58%  '
